$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instructions")

# Update Instruction Type values for rows 3-5 ("R1 R2" -> "R1 R2 B")
$ws.Range("E3").Value = "R1 R2 B"
$ws.Range("E4").Value = "R1 R2 B"
$ws.Range("E5").Value = "R1 R2 B"

# Update Instruction Type value for row 15 ("R1 O" -> "R1")
$ws.Range("E15").Value = "R1"

# Update the active sheet view / selection
$ws.Activate()
$ws.Range("E16").Select()
